# Auto update ESO load data
# Updates existing load rows (2026-02-07 .. 2026-02-11) with refreshed
# forecast values and appends a new row for 2026-02-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 2026
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = 7
$ws.Cells.Item(4, 4).Value = 4778
$ws.Cells.Item(4, 5).Value = 4496
$ws.Cells.Item(4, 6).Value = 4379
$ws.Cells.Item(4, 7).Value = 4393
$ws.Cells.Item(4, 8).Value = 4389
$ws.Cells.Item(4, 9).Value = 4412
$ws.Cells.Item(4, 10).Value = 4555
$ws.Cells.Item(4, 11).Value = 4947
$ws.Cells.Item(4, 12).Value = 5381
$ws.Cells.Item(4, 13).Value = 5702
$ws.Cells.Item(4, 14).Value = 5785
$ws.Cells.Item(4, 15).Value = 5741
$ws.Cells.Item(4, 16).Value = 5709
$ws.Cells.Item(4, 17).Value = 5631
$ws.Cells.Item(4, 18).Value = 5566
$ws.Cells.Item(4, 19).Value = 5561
$ws.Cells.Item(4, 20).Value = 5605
$ws.Cells.Item(4, 21).Value = 5694
$ws.Cells.Item(4, 22).Value = 5769
$ws.Cells.Item(4, 23).Value = 5646
$ws.Cells.Item(4, 24).Value = 5386
$ws.Cells.Item(4, 25).Value = 5108
$ws.Cells.Item(4, 26).Value = 4971
$ws.Cells.Item(4, 27).Value = 4776
$ws.Cells.Item(5, 1).Value = 2026
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 8
$ws.Cells.Item(5, 4).Value = 4646
$ws.Cells.Item(5, 5).Value = 4393
$ws.Cells.Item(5, 6).Value = 4276
$ws.Cells.Item(5, 7).Value = 4255
$ws.Cells.Item(5, 8).Value = 4318
$ws.Cells.Item(5, 9).Value = 4337
$ws.Cells.Item(5, 10).Value = 4442
$ws.Cells.Item(5, 11).Value = 4670
$ws.Cells.Item(5, 12).Value = 5022
$ws.Cells.Item(5, 13).Value = 5362
$ws.Cells.Item(5, 14).Value = 5567
$ws.Cells.Item(5, 15).Value = 5661
$ws.Cells.Item(5, 16).Value = 5711
$ws.Cells.Item(5, 17).Value = 5562
$ws.Cells.Item(5, 18).Value = 5562
$ws.Cells.Item(5, 19).Value = 5577
$ws.Cells.Item(5, 20).Value = 5646
$ws.Cells.Item(5, 21).Value = 5776
$ws.Cells.Item(5, 22).Value = 5911
$ws.Cells.Item(5, 23).Value = 5791
$ws.Cells.Item(5, 24).Value = 5547
$ws.Cells.Item(5, 25).Value = 5272
$ws.Cells.Item(5, 26).Value = 5071
$ws.Cells.Item(5, 27).Value = 4789
$ws.Cells.Item(6, 1).Value = 2026
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = 9
$ws.Cells.Item(6, 4).Value = 4926
$ws.Cells.Item(6, 5).Value = 4675
$ws.Cells.Item(6, 6).Value = 4565
$ws.Cells.Item(6, 7).Value = 4578
$ws.Cells.Item(6, 8).Value = 4607
$ws.Cells.Item(6, 9).Value = 4731
$ws.Cells.Item(6, 10).Value = 5121
$ws.Cells.Item(6, 11).Value = 5631
$ws.Cells.Item(6, 12).Value = 6130
$ws.Cells.Item(6, 13).Value = 6394
$ws.Cells.Item(6, 14).Value = 6419
$ws.Cells.Item(6, 15).Value = 6359
$ws.Cells.Item(6, 16).Value = 6274
$ws.Cells.Item(6, 17).Value = 6249
$ws.Cells.Item(6, 18).Value = 6199
$ws.Cells.Item(6, 19).Value = 6139
$ws.Cells.Item(6, 20).Value = 6072
$ws.Cells.Item(6, 21).Value = 6167
$ws.Cells.Item(6, 22).Value = 6362
$ws.Cells.Item(6, 23).Value = 6336
$ws.Cells.Item(6, 24).Value = 6120
$ws.Cells.Item(6, 25).Value = 5816
$ws.Cells.Item(6, 26).Value = 5534
$ws.Cells.Item(6, 27).Value = 5199
$ws.Cells.Item(7, 1).Value = 2026
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = 10
$ws.Cells.Item(7, 4).Value = 5304
$ws.Cells.Item(7, 5).Value = 5066
$ws.Cells.Item(7, 6).Value = 4988
$ws.Cells.Item(7, 7).Value = 4962
$ws.Cells.Item(7, 8).Value = 4972
$ws.Cells.Item(7, 9).Value = 5135
$ws.Cells.Item(7, 10).Value = 5597
$ws.Cells.Item(7, 11).Value = 6134
$ws.Cells.Item(7, 12).Value = 6593
$ws.Cells.Item(7, 13).Value = 6745
$ws.Cells.Item(7, 14).Value = 6671
$ws.Cells.Item(7, 15).Value = 6557
$ws.Cells.Item(7, 16).Value = 6444
$ws.Cells.Item(7, 17).Value = 6385
$ws.Cells.Item(7, 18).Value = 6306
$ws.Cells.Item(7, 19).Value = 6328
$ws.Cells.Item(7, 20).Value = 6363
$ws.Cells.Item(7, 21).Value = 6481
$ws.Cells.Item(7, 22).Value = 6678
$ws.Cells.Item(7, 23).Value = 6647
$ws.Cells.Item(7, 24).Value = 6424
$ws.Cells.Item(7, 25).Value = 6117
$ws.Cells.Item(7, 26).Value = 5852
$ws.Cells.Item(7, 27).Value = 5534
$ws.Cells.Item(8, 1).Value = 2026
$ws.Cells.Item(8, 2).Value = 2
$ws.Cells.Item(8, 3).Value = 11
$ws.Cells.Item(8, 4).Value = 5311
$ws.Cells.Item(8, 5).Value = 5073
$ws.Cells.Item(8, 6).Value = 4994
$ws.Cells.Item(8, 7).Value = 4969
$ws.Cells.Item(8, 8).Value = 4979
$ws.Cells.Item(8, 9).Value = 5141
$ws.Cells.Item(8, 10).Value = 5604
$ws.Cells.Item(8, 11).Value = 6141
$ws.Cells.Item(8, 12).Value = 6600
$ws.Cells.Item(8, 13).Value = 6752
$ws.Cells.Item(8, 14).Value = 6677
$ws.Cells.Item(8, 15).Value = 6563
$ws.Cells.Item(8, 16).Value = 6450
$ws.Cells.Item(8, 17).Value = 6390
$ws.Cells.Item(8, 18).Value = 6311
$ws.Cells.Item(8, 19).Value = 6333
$ws.Cells.Item(8, 20).Value = 6368
$ws.Cells.Item(8, 21).Value = 6487
$ws.Cells.Item(8, 22).Value = 6685
$ws.Cells.Item(8, 23).Value = 6654
$ws.Cells.Item(8, 24).Value = 6430
$ws.Cells.Item(8, 25).Value = 6124
$ws.Cells.Item(8, 26).Value = 5859
$ws.Cells.Item(8, 27).Value = 5541
$ws.Cells.Item(9, 1).Value = 2026
$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = 12
$ws.Cells.Item(9, 4).Value = 5121
$ws.Cells.Item(9, 5).Value = 4884
$ws.Cells.Item(9, 6).Value = 4806
$ws.Cells.Item(9, 7).Value = 4780
$ws.Cells.Item(9, 8).Value = 4790
$ws.Cells.Item(9, 9).Value = 4953
$ws.Cells.Item(9, 10).Value = 5414
$ws.Cells.Item(9, 11).Value = 5950
$ws.Cells.Item(9, 12).Value = 6408
$ws.Cells.Item(9, 13).Value = 6559
$ws.Cells.Item(9, 14).Value = 6494
$ws.Cells.Item(9, 15).Value = 6394
$ws.Cells.Item(9, 16).Value = 6295
$ws.Cells.Item(9, 17).Value = 6242
$ws.Cells.Item(9, 18).Value = 6173
$ws.Cells.Item(9, 19).Value = 6192
$ws.Cells.Item(9, 20).Value = 6223
$ws.Cells.Item(9, 21).Value = 6326
$ws.Cells.Item(9, 22).Value = 6498
$ws.Cells.Item(9, 23).Value = 6467
$ws.Cells.Item(9, 24).Value = 6241
$ws.Cells.Item(9, 25).Value = 5932
$ws.Cells.Item(9, 26).Value = 5665
$ws.Cells.Item(9, 27).Value = 5343
